$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 33; this shifts existing rows 33-118 down to 34-119,
# carrying their formatting and values with them (Excel's native row-insert shift
# behavior), and automatically extends the sheet's used range / dimension.
$ws.Rows("33:33").Insert()

# Populate the newly-inserted (blank) row 33 with the new weekly record.
$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44497
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = 100112029
$ws.Range("G33").Value = "Orégano"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 36
$ws.Range("K33").Value = 8000
$ws.Range("L33").Value = 9000
$ws.Range("M33").Value = 8472
$ws.Range("N33").Value = "`$/docena de atados"
$ws.Range("O33").Value = "Región Metropolitana"
$ws.Range("P33").Value = 2824
$ws.Range("Q33").Value = 3
$ws.Range("R33").Value = "Hortaliza"

# Ensure the date cell keeps the workbook's date number format (style s="2"),
# matching every other row in column D.
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
